# Add a new "The god father" row (row 13) to the cost sheet, mirroring the
# layout/formatting of the row directly above it (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting (style) of row 12's populated cells down to row 13 ---
# Doing this cell-by-cell (rather than whole-row) avoids manufacturing spurious
# empty cells in columns that row 12 leaves blank (F, G, H, J, K, L..Q).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("E12").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null

$ws.Range("I12").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4122) | Out-Null

# --- Fill in the new row's content ---
$ws.Range("A13").Value = "The god father"
$ws.Range("B13").Value = 17
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Formula = "=B13*500"
$ws.Range("I13").Formula = "=B13*200"

# --- Match the author's final cursor/selection position ---
$ws.Range("D13").Select() | Out-Null
